# Fixed M_i(L) to N_i(L)
#
# Column A ("L") values for rows 36-101 had lost their "tens"/"twenties"
# component (a sawtooth reset every ~10 rows). Restore the correct
# continuous sequence: rows 36-68 belong to the "10s" decade (+10), rows
# 69-101 belong to the "20s" decade (+20). Rows that already read
# correctly (the decade anchors 12, 15, 18, 21, 24, 27) are left
# untouched, matching the source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 36; $r -le 101; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value2
    if ($v -ne $null -and $v -lt 10) {
        if ($r -le 68) {
            $cell.Value = $v + 10
        } else {
            $cell.Value = $v + 20
        }
    }
}

# Column A got a touch wider after the longer two-digit values landed
# (bestFit recalculates once the two-digit numbers push past the old width).
$ws.Columns.Item(1).ColumnWidth = 4.25

# Selection moved from C2 to O6.
$ws.Range("O6").Select()
